$wb = $excel.ActiveWorkbook

# The workbook already holds several identical "severalAnswersN" sheets
# (same Identifiers/grid layout). Add one more of the same kind,
# "severalAnswers4", at the end of the tab strip.
$template = $wb.Worksheets.Item("severalAnswers2")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $last)

$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "severalAnswers4"
